$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '33.996.23'
$ws.Range("E2").Value = '  -0.15%  '

$ws.Range("D3").Value = '1.777.87'
$ws.Range("E3").Value = '  -2.13%  '

$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '225.89'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.54%  '

$ws.Range("E6").Value = '  -0.94%  '

$ws.Range("E7").Value = '  -0.02%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '32.42'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +2.05%  '

$ws.Range("E9").Value = '  -1.49%  '

$ws.Range("E10").Value = '  -2.54%  '

$ws.Range("E11").Value = '  +0.68%  '

$ws.Range("D12").Value = '2.035.82'
$ws.Range("E12").Value = '  -2.52%  '

$ws.Range("D13").Value = '1.782.45'
$ws.Range("E13").Value = '  -2.69%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '10.87'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +0.93%  '

$ws.Range("D15").Value = '33.997.94'
$ws.Range("E15").Value = '  -0.41%  '

$ws.Range("E16").Value = '  -3.71%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '4.11'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -4.19%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '67.60'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -2.36%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '243.30'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -2.36%  '

$ws.Range("E20").Value = '  -0.80%  '

$ws.Range("E21").Value = '  +0.07%  '

$ws.Range("E22").Value = '  -3.23%  '

$ws.Range("E23").Value = '  -4.13%  '

$ws.Range("E24").Value = '  -3.51%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '159.87'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -0.16%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '16.24'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -1.42%  '

$ws.Range("E27").Value = '  -2.45%  '

$ws.Range("E28").Value = '  -1.96%  '

$ws.Range("E29").Value = '  -0.03%  '

$ws.Range("E30").Value = '  +0.82%  '

$ws.Range("E31").Value = '  -3.56%  '

$ws.Range("E32").Value = '  -3.59%  '

$ws.Range("E33").Value = '  -2.03%  '

$ws.Range("E34").Value = '  -4.27%  '

$ws.Range("D35").Value = '1.384.21'
$ws.Range("E35").Value = '  -3.18%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.644'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +1.43%  '

$ws.Range("E37").Value = '  -1.69%  '

$ws.Range("E38").Value = '  -1.54%  '

$ws.Range("E39").Value = '  +0.07%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.18'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +1.77%  '

$ws.Range("B41").Value = 'ARBITRUM'
$ws.Range("C41").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.908'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -4.83%  '

$ws.Range("B42").Value = 'MXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.68'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -4.47%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '77.33'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -4.58%  '

$ws.Range("E44").Value = '  +14.58%  '

$ws.Range("E45").Value = '  +2.22%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.51'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +5.96%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0498'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +0.25%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '107.31'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +0.72%  '

$ws.Range("E49").Value = '  -4.22%  '

$ws.Range("D50").Value = '1.935.00'
$ws.Range("E50").Value = '  -2.48%  '

$ws.Range("E51").Value = '  +0.13%  '
